$wb = $excel.ActiveWorkbook

# Rename the "Device" sheet to "Apparatus"
$ws = $wb.Worksheets.Item("Device")
$ws.Name = "Apparatus"

# Update the summary sentence on the renamed sheet
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Make the header row (bus number / type number / parameters) bold,
# matching the style already used elsewhere (copy format from A2, which
# already carries that bold style), instead of inventing a new style.
$ws.Range("A2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)

# Make "Apparatus" the active sheet/tab with A2 selected
$null = $ws.Activate()
$null = $ws.Range("A2").Select()
